$wb = $excel.ActiveWorkbook

# --- Sheet: Teacher ---
$ws1 = $wb.Worksheets.Item("Teacher")
$ws1.Range("A2").Value = 4
$ws1.Range("B2").Value = "Zendaya"
$ws1.Range("C2").Value = "Guru kelas 5A"
$ws1.Range("A3").Value = 2
$ws1.Range("B3").Value = "Tom Holland"
$ws1.Range("C3").Value = "Guru kelas 4A"
$ws1.Range("A4").Value = 3
$ws1.Range("B4").Value = "Andrew Garfield"
$ws1.Range("C4").Value = "Guru kelas 4B"
$ws1.Range("A5").Value = 5
$ws1.Range("B5").Value = "Crystenz Danz"
$ws1.Range("C5").Value = "Guru kelas 5B"
$ws1.Range("A6").Value = 1
$ws1.Range("B6").Value = "Emma Stone"
$ws1.Range("C6").Value = "Guru kelas 3A dan 3B"

# --- Sheet: StudentAssessment ---
$ws2 = $wb.Worksheets.Item("StudentAssessment")
$ws2.Range("A2").Value = 2
$ws2.Range("B2").Value = "Siti Nurhaliza"
$ws2.Range("C2").Value = "3A"
$ws2.Range("D2").Value = 1
$ws2.Range("E2").Value = 1
$ws2.Range("F2").Value = 5
$ws2.Range("G2").Value = 5
$ws2.Range("H2").Value = 4
$ws2.Range("I2").Value = 4
$ws2.Range("J2").Value = 4
$ws2.Range("K2").Value = 22
$ws2.Range("L2").Value = 5
$ws2.Range("M2").Value = 5
$ws2.Range("N2").Value = 5
$ws2.Range("O2").Value = 4
$ws2.Range("P2").Value = 5
$ws2.Range("Q2").Value = 24
$ws2.Range("R2").Value = 5
$ws2.Range("S2").Value = 5
$ws2.Range("T2").Value = 5
$ws2.Range("U2").Value = 5
$ws2.Range("V2").Value = 5
$ws2.Range("W2").Value = 25
$ws2.Range("X2").Value = 71
$ws2.Range("Y2").Value = 23.67
$ws2.Range("Z2").Value = "Sangat Baik"
$ws2.Range("AA2").Value = "Siswa sangat aktif dan menunjukkan kemampuan yang sangat baik di semua aspek."
$ws2.Range("AB2").Value = 45907.19290165509
$ws2.Range("AC2").Value = 45907.19290165509
$ws2.Range("A3").Value = 3
$ws2.Range("B3").Value = "Budi Santoso"
$ws2.Range("C3").Value = "3B"
$ws2.Range("D3").Value = 1
$ws2.Range("E3").Value = 1
$ws2.Range("F3").Value = 4
$ws2.Range("G3").Value = 3
$ws2.Range("H3").Value = 3
$ws2.Range("I3").Value = 3
$ws2.Range("J3").Value = 3
$ws2.Range("K3").Value = 16
$ws2.Range("L3").Value = 4
$ws2.Range("M3").Value = 4
$ws2.Range("N3").Value = 3
$ws2.Range("O3").Value = 3
$ws2.Range("P3").Value = 4
$ws2.Range("Q3").Value = 18
$ws2.Range("R3").Value = 5
$ws2.Range("S3").Value = 4
$ws2.Range("T3").Value = 4
$ws2.Range("U3").Value = 4
$ws2.Range("V3").Value = 4
$ws2.Range("W3").Value = 21
$ws2.Range("X3").Value = 55
$ws2.Range("Y3").Value = 18.33
$ws2.Range("Z3").Value = "Baik"
$ws2.Range("AA3").Value = "Siswa menunjukkan peningkatan bertahap, perlu lebih banyak latihan pengucapan."
$ws2.Range("AB3").Value = 45907.19290210648
$ws2.Range("AC3").Value = 45907.19290210648
$ws2.Range("A4").Value = 8
$ws2.Range("B4").Value = "Fitri Ramadhani"
$ws2.Range("C4").Value = "5A"
$ws2.Range("D4").Value = 4
$ws2.Range("E4").Value = 1
$ws2.Range("F4").Value = 5
$ws2.Range("G4").Value = 3
$ws2.Range("H4").Value = 4
$ws2.Range("I4").Value = 3
$ws2.Range("J4").Value = 3
$ws2.Range("K4").Value = 18
$ws2.Range("L4").Value = 5
$ws2.Range("M4").Value = 4
$ws2.Range("N4").Value = 4
$ws2.Range("O4").Value = 4
$ws2.Range("P4").Value = 4
$ws2.Range("Q4").Value = 21
$ws2.Range("R4").Value = 5
$ws2.Range("S4").Value = 4
$ws2.Range("T4").Value = 4
$ws2.Range("U4").Value = 4
$ws2.Range("V4").Value = 4
$ws2.Range("W4").Value = 21
$ws2.Range("X4").Value = 60
$ws2.Range("Y4").Value = 20
$ws2.Range("Z4").Value = "Baik"
$ws2.Range("AA4").Value = "Siswa rajin hadir dan menunjukkan peningkatan yang konsisten."
$ws2.Range("AB4").Value = 45907.19290434028
$ws2.Range("AC4").Value = 45907.19290434028
$ws2.Range("A5").Value = 9
$ws2.Range("B5").Value = "Joko Susilo"
$ws2.Range("C5").Value = "5B"
$ws2.Range("D5").Value = 5
$ws2.Range("E5").Value = 1
$ws2.Range("F5").Value = 3
$ws2.Range("G5").Value = 2
$ws2.Range("H5").Value = 2
$ws2.Range("I5").Value = 2
$ws2.Range("J5").Value = 2
$ws2.Range("K5").Value = 11
$ws2.Range("L5").Value = 3
$ws2.Range("M5").Value = 3
$ws2.Range("N5").Value = 2
$ws2.Range("O5").Value = 2
$ws2.Range("P5").Value = 3
$ws2.Range("Q5").Value = 13
$ws2.Range("R5").Value = 4
$ws2.Range("S5").Value = 3
$ws2.Range("T5").Value = 3
$ws2.Range("U5").Value = 3
$ws2.Range("V5").Value = 3
$ws2.Range("W5").Value = 16
$ws2.Range("X5").Value = 40
$ws2.Range("Y5").Value = 13.33
$ws2.Range("Z5").Value = "Cukup"
$ws2.Range("AA5").Value = "Siswa memerlukan perhatian khusus dan motivasi tambahan untuk meningkatkan partisipasi."
$ws2.Range("AB5").Value = 45907.19290532407
$ws2.Range("AC5").Value = 45907.19290532407
$ws2.Range("A6").Value = 10
$ws2.Range("B6").Value = "Indah Permata"
$ws2.Range("C6").Value = "5B"
$ws2.Range("D6").Value = 5
$ws2.Range("E6").Value = 1
$ws2.Range("F6").Value = 5
$ws2.Range("G6").Value = 5
$ws2.Range("H6").Value = 5
$ws2.Range("I6").Value = 4
$ws2.Range("J6").Value = 5
$ws2.Range("K6").Value = 24
$ws2.Range("L6").Value = 5
$ws2.Range("M6").Value = 5
$ws2.Range("N6").Value = 5
$ws2.Range("O6").Value = 5
$ws2.Range("P6").Value = 5
$ws2.Range("Q6").Value = 25
$ws2.Range("R6").Value = 5
$ws2.Range("S6").Value = 5
$ws2.Range("T6").Value = 5
$ws2.Range("U6").Value = 5
$ws2.Range("V6").Value = 5
$ws2.Range("W6").Value = 25
$ws2.Range("X6").Value = 74
$ws2.Range("Y6").Value = 24.67
$ws2.Range("Z6").Value = "Sangat Baik"
$ws2.Range("AA6").Value = "Siswa sangat berprestasi dan dapat membantu teman-teman yang kesulitan."
$ws2.Range("AB6").Value = 45907.19290578704
$ws2.Range("AC6").Value = 45907.19290578704
$ws2.Range("A7").Value = 5
$ws2.Range("B7").Value = "Andi Wijaya"
$ws2.Range("C7").Value = "4A"
$ws2.Range("D7").Value = 2
$ws2.Range("E7").Value = 1
$ws2.Range("F7").Value = 3
$ws2.Range("G7").Value = 3
$ws2.Range("H7").Value = 2
$ws2.Range("I7").Value = 2
$ws2.Range("J7").Value = 2
$ws2.Range("K7").Value = 12
$ws2.Range("L7").Value = 4
$ws2.Range("M7").Value = 3
$ws2.Range("N7").Value = 3
$ws2.Range("O7").Value = 3
$ws2.Range("P7").Value = 3
$ws2.Range("Q7").Value = 16
$ws2.Range("R7").Value = 4
$ws2.Range("S7").Value = 4
$ws2.Range("T7").Value = 3
$ws2.Range("U7").Value = 3
$ws2.Range("V7").Value = 3
$ws2.Range("W7").Value = 17
$ws2.Range("X7").Value = 45
$ws2.Range("Y7").Value = 15
$ws2.Range("Z7").Value = "Cukup"
$ws2.Range("AA7").ClearContents()
$ws2.Range("AB7").Value = 45907.19290296296
$ws2.Range("AC7").Value = 45907.20727947917
$ws2.Range("A8").Value = 1
$ws2.Range("B8").Value = "Ahmad Rizki"
$ws2.Range("C8").Value = "3A"
$ws2.Range("D8").Value = 1
$ws2.Range("E8").Value = 1
$ws2.Range("F8").Value = 5
$ws2.Range("G8").Value = 4
$ws2.Range("H8").Value = 4
$ws2.Range("I8").Value = 3
$ws2.Range("J8").Value = 4
$ws2.Range("K8").Value = 20
$ws2.Range("L8").Value = 5
$ws2.Range("M8").Value = 4
$ws2.Range("N8").Value = 5
$ws2.Range("O8").Value = 4
$ws2.Range("P8").Value = 4
$ws2.Range("Q8").Value = 22
$ws2.Range("R8").Value = 5
$ws2.Range("S8").Value = 5
$ws2.Range("T8").Value = 4
$ws2.Range("U8").Value = 4
$ws2.Range("V8").Value = 2
$ws2.Range("W8").Value = 20
$ws2.Range("X8").Value = 62
$ws2.Range("Y8").Value = 20.67
$ws2.Range("Z8").Value = "Baik"
$ws2.Range("AA8").ClearContents()
$ws2.Range("AB8").Value = 45907.192900625
$ws2.Range("AC8").Value = 45907.24896291667
$ws2.Range("A9").Value = 7
$ws2.Range("B9").Value = "Rudi Hermawan"
$ws2.Range("C9").Value = "5A"
$ws2.Range("D9").Value = 4
$ws2.Range("E9").Value = 1
$ws2.Range("F9").Value = 5
$ws2.Range("G9").Value = 5
$ws2.Range("H9").Value = 5
$ws2.Range("I9").Value = 5
$ws2.Range("J9").Value = 5
$ws2.Range("K9").Value = 25
$ws2.Range("L9").Value = 5
$ws2.Range("M9").Value = 5
$ws2.Range("N9").Value = 5
$ws2.Range("O9").Value = 5
$ws2.Range("P9").Value = 5
$ws2.Range("Q9").Value = 25
$ws2.Range("R9").Value = 0
$ws2.Range("S9").Value = 0
$ws2.Range("T9").Value = 0
$ws2.Range("U9").Value = 0
$ws2.Range("V9").Value = 0
$ws2.Range("W9").Value = 0
$ws2.Range("X9").Value = 50
$ws2.Range("Y9").Value = 16.67
$ws2.Range("Z9").Value = "Baik"
$ws2.Range("AA9").ClearContents()
$ws2.Range("AB9").Value = 45907.19290392361
$ws2.Range("AC9").Value = 45907.25062901621
$ws2.Range("A10").Value = 6
$ws2.Range("B10").Value = "Maya Sari"
$ws2.Range("C10").Value = "4B"
$ws2.Range("D10").Value = 3
$ws2.Range("E10").Value = 1
$ws2.Range("F10").Value = 5
$ws2.Range("G10").Value = 4
$ws2.Range("H10").Value = 4
$ws2.Range("I10").Value = 4
$ws2.Range("J10").Value = 4
$ws2.Range("K10").Value = 21
$ws2.Range("L10").Value = 5
$ws2.Range("M10").Value = 4
$ws2.Range("N10").Value = 4
$ws2.Range("O10").Value = 4
$ws2.Range("P10").Value = 4
$ws2.Range("Q10").Value = 21
$ws2.Range("R10").Value = 5
$ws2.Range("S10").Value = 5
$ws2.Range("T10").Value = 4
$ws2.Range("U10").Value = 4
$ws2.Range("V10").Value = 5
$ws2.Range("W10").Value = 23
$ws2.Range("X10").Value = 65
$ws2.Range("Y10").Value = 21.67
$ws2.Range("Z10").Value = "Sangat Baik"
$ws2.Range("AA10").ClearContents()
$ws2.Range("AB10").Value = 45907.19290341435
$ws2.Range("AC10").Value = 45907.25089626158
$ws2.Range("A11").Value = 12
$ws2.Range("B11").Value = "tahap 2"
$ws2.Range("C11").Value = "5A"
$ws2.Range("D11").Value = 4
$ws2.Range("E11").Value = 6
$ws2.Range("F11").Value = 5
$ws2.Range("G11").Value = 5
$ws2.Range("H11").Value = 5
$ws2.Range("I11").Value = 5
$ws2.Range("J11").Value = 5
$ws2.Range("K11").Value = 25
$ws2.Range("L11").Value = 5
$ws2.Range("M11").Value = 5
$ws2.Range("N11").Value = 5
$ws2.Range("O11").Value = 5
$ws2.Range("P11").Value = 5
$ws2.Range("Q11").Value = 25
$ws2.Range("R11").Value = 0
$ws2.Range("S11").Value = 0
$ws2.Range("T11").Value = 0
$ws2.Range("U11").Value = 0
$ws2.Range("V11").Value = 0
$ws2.Range("W11").Value = 0
$ws2.Range("X11").Value = 50
$ws2.Range("Y11").Value = 25
$ws2.Range("Z11").Value = "Sangat Baik"
$ws2.Range("AA11").Value = "tahap2"
$ws2.Range("AB11").Value = 45911.25347222222
$ws2.Range("AC11").Value = 45907.32526325231
$ws2.Range("A12").Value = 4
$ws2.Range("B12").Value = "Dewi Lestari"
$ws2.Range("C12").Value = "4A"
$ws2.Range("D12").Value = 2
$ws2.Range("E12").Value = 1
$ws2.Range("F12").Value = 5
$ws2.Range("G12").Value = 5
$ws2.Range("H12").Value = 5
$ws2.Range("I12").Value = 5
$ws2.Range("J12").Value = 5
$ws2.Range("K12").Value = 25
$ws2.Range("L12").Value = 5
$ws2.Range("M12").Value = 5
$ws2.Range("N12").Value = 5
$ws2.Range("O12").Value = 5
$ws2.Range("P12").Value = 5
$ws2.Range("Q12").Value = 25
$ws2.Range("R12").Value = 5
$ws2.Range("S12").Value = 5
$ws2.Range("T12").Value = 5
$ws2.Range("U12").Value = 5
$ws2.Range("V12").Value = 4
$ws2.Range("W12").Value = 24
$ws2.Range("X12").Value = 74
$ws2.Range("Y12").Value = 24.67
$ws2.Range("Z12").Value = "Sangat Baik"
$ws2.Range("AA12").Value = "Siswa sangat berprestasi dan konsisten dalam semua aspek pembelajaran."
$ws2.Range("AB12").Value = 45907.19290255787
$ws2.Range("AC12").Value = 45907.26760246528

# --- Sheet: _prisma_migrations ---
$ws3 = $wb.Worksheets.Item("_prisma_migrations")
$ws3.Range("A2").Value = "965e5a3d-e2e7-4986-930b-85e67d978f96"
$ws3.Range("C2").Value = 45907.19287397835
$ws3.Range("G2").Value = 45907.19286892472
